$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 42
$prevRow = $row - 1

# New log entry values (columns A-H)
$ws.Cells.Item($row, 1).Value = "2025-08-22 03:51:15 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-22 09:21:15 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Match the formatting of the preceding data row
$ws.Range("A" + $prevRow + ":H" + $prevRow).Copy()
$ws.Range("A" + $row + ":H" + $row).PasteSpecial(-4122)
